$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 'Paula Belén Chairez Rosas_20251120_205610'

$c = $ws.Range("B13")
$c.Formula = "'"
$c.Style = "Normal"

$ws.Range("C13").Value = 'Paula Belén Chairez Rosas'

$ws.Range("D13").Value = 20

$ws.Range("E13").Value = 'Femenino'

$ws.Range("F13").Value = '2025-11-20 20:56:10'

$ws.Range("G13").Value = '{
  "portion": 0.6,
  "diet": 0.42857142857142855,
  "salt": 0.2,
  "fat": 0.4,
  "natural": 0.6,
  "convenience": 0.8,
  "price": 0.8
}'

$ws.Range("H13").Value = 'Nongshim Neoguri Spicy Seafood'

$c = $ws.Range("I13")
$c.NumberFormat = "@"
$c.Value = '0.637'
$c.Style = "Normal"

$ws.Range("J13").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'

$ws.Range("K13").Value = 'Nissin Chow Mein Teriyaki Beef'

$c = $ws.Range("L13")
$c.NumberFormat = "@"
$c.Value = '0.522'
$c.Style = "Normal"

$ws.Range("M13").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'

$ws.Range("N13").Value = 'Maruchan Ramen Sabor Pollo'

$c = $ws.Range("O13")
$c.NumberFormat = "@"
$c.Value = '0.515'
$c.Style = "Normal"

$ws.Range("P13").Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'

$ws.Range("Q13").Value = 'Kraft Macaroni & Cheese Dinner'

$c = $ws.Range("R13")
$c.NumberFormat = "@"
$c.Value = '0.556'
$c.Style = "Normal"

$ws.Range("S13").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'

$ws.Range("T13").Value = 'Amy’s Macaroni & Cheese (frozen)'

$c = $ws.Range("U13")
$c.NumberFormat = "@"
$c.Value = '0.548'
$c.Style = "Normal"

$ws.Range("V13").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'

$ws.Range("W13").Value = 'Velveeta Original Shells & Cheese (microwave cups)'

$c = $ws.Range("X13")
$c.NumberFormat = "@"
$c.Value = '0.521'
$c.Style = "Normal"

$ws.Range("Y13").Value = 'Muy cremoso, porción individual, rápido, salado, ideal para niños'

$ws.Range("Z13").Value = 'Wild Planet Wild Tuna Pasta Salad'

$c = $ws.Range("AA13")
$c.NumberFormat = "@"
$c.Value = '0.684'
$c.Style = "Normal"

$ws.Range("AB13").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'

$ws.Range("AC13").Value = 'StarKist Chicken Creations (Chicken Salad)'

$c = $ws.Range("AD13")
$c.NumberFormat = "@"
$c.Value = '0.582'
$c.Style = "Normal"

$ws.Range("AE13").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'

$ws.Range("AF13").Value = 'Jack Link’s Beef Jerky Original'

$c = $ws.Range("AG13")
$c.NumberFormat = "@"
$c.Value = '0.575'
$c.Style = "Normal"

$ws.Range("AH13").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

$ws.Rows(13).EntireRow.AutoFit()
